# 1st iteration (#141): update the generated ValueSet metadata sheet.
#  - "Experimental" row: set the Value cell (B7) to the text "true"
#  - "Date" row: update the Value cell (B8) to the new generation timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Writing the literal string "true" straight into Value/Value2 gets smart-
# converted to a Boolean by Excel (same as typing TRUE into a cell), so we
# build it as text via a formula result on a scratch cell and paste just the
# value back in - this keeps the cell's type as text ("true") without
# disturbing its existing cell style.
$scratch = $ws.Range("Z100")
$scratch.Formula = '="true"'
$scratch.Copy()
$ws.Range("B7").PasteSpecial(-4163)
$scratch.Clear()

$ws.Range("B8").Value = "2025-07-21T12:46:15+00:00"
